# Auto-generated Excel COM-interop script
# Updates Price (D) and Volume(1h) (E) columns for the crypto price tracker sheet
# to match the new snapshot of values scraped on 2023-01-20.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Using NumberFormat "@" (Text) before
# assignment keeps these values as literal strings (matching the original
# inline-string cells) instead of letting Excel auto-convert them to numbers
# or percentages; resetting the style back to "Normal" afterwards keeps the
# cell formatting unchanged.
$updates = @{
    "D2" = "289.42"
    "E2" = "-0.17%"
    "D3" = "31.09"
    "E3" = "1.74%"
    "D4" = "4.946"
    "E4" = "0.58%"
    "D5" = "0.07362"
    "E5" = "1.54%"
    "D6" = "2.309"
    "E6" = "27.54%"
    "D7" = "7.670"
    "E7" = "0.45%"
    "D8" = "0.9183"
    "E8" = "1.93%"
    "D9" = "0.09149"
    "E9" = "13.96%"
    "D10" = "0.1701"
    "E10" = "0.79%"
    "D11" = "0.08271"
    "E11" = "1.54%"
    "D12" = "0.03112"
    "E12" = "1.92%"
    "D13" = "0.09992"
    "E13" = "-0.13%"
    "D14" = "0.001496"
    "E14" = "-0.01%"
    "D15" = "0.005748"
    "E15" = "0.85%"
    "D16" = "3.472"
    "D17" = "3.741"
    "E17" = "1.07%"
    "D18" = "2.103"
    "E18" = "1.35%"
    "D19" = "0.3332"
    "E19" = "0.39%"
    "D20" = "0.1300"
    "E20" = "-0.18%"
    "D21" = "4.179"
    "E21" = "5.35%"
    "D22" = "0.2123"
    "E22" = "-2.10%"
    "D23" = "0.04509"
    "E23" = "0.06%"
    "D24" = "0.001215"
    "E24" = "0.29%"
    "D25" = "0.004196"
    "E25" = "-5.42%"
    "D26" = "0.0001300"
    "E26" = "0.01%"
    "D27" = "0.0003393"
    "E27" = "0.16%"
    "D39" = "0.01573"
    "E39" = "-0.66%"
    "D40" = "0.04508"
    "E40" = "3.57%"
    "D41" = "0.007366"
    "E41" = "1.22%"
    "D42" = "0.009849"
    "E42" = "-1.74%"
    "D43" = "0.1337"
    "E43" = "1.73%"
    "D44" = "0.002220"
    "E44" = "10.78%"
    "D45" = "0.008516"
    "E45" = "-9.88%"
    "E46" = "4.22%"
    "D47" = "0.00000000750"
    "E47" = "0.17%"
    "D48" = "2.445"
    "E48" = "8.46%"
    "D49" = "0.001999"
    "E49" = "-30.89%"
    "D50" = "0.00002099"
    "E50" = "0.17%"
    "D51" = "0.0001999"
    "E51" = "0.17%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

